$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 20 (everything from old row 20 downward shifts
# down by one: old r20->r21 ... old r31->r32), mirroring Excel's
# "insert row" behaviour so existing formatting below stays intact.
$ws.Rows("20:20").Insert()

# Carry over the row 19 formatting (cell styles + row height) onto the new
# row 20 before filling in its own values, so s="7"/"5"/"8" + ht are right.
$ws.Range("A19:C19").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# New resource entry: S[14] / Informix / Software
$ws.Range("A20").Value = "S[14]"
$ws.Range("B20").Value = "Informix"
$ws.Range("C20").Value = "Software"

# Match the row height used by every other data row in the table.
$ws.Rows.Item(20).RowHeight = 18.75

# Restore the view: scroll back to the top and leave the selection on B13,
# matching the saved workbook state.
[void]$ws.Range("B13").Select()
